{"js": "// Slight change to samtools exercises\n// 1) \"samtools view -bS sample1.sam > sample1.bam\"\n//       -> \"samtools view -S -b sample1.sam > sample1.bam\"\n// 2) \"samtools sort sample1.bam -o sample1.bam # sort bam file\"\n//       -> \"samtools sort sample1.bam sample1 # sort bam file\"\n\nconst body = context.document.body;\n\n// --- Edit 1: \"-bS \" -> \"-S -b \" on the samtools view line ---\nconst viewHits = body.search(\"-bS \", { matchCase: true, matchWholeWord: false });\nviewHits.load(\"text\");\nawait context.sync();\n\nif (viewHits.items.length > 0) {\n  viewHits.items[0].insertText(\"-S -b \", \"Replace\");\n  await context.sync();\n}\n\n// --- Edit 2: \"sort sample1.bam -o sample1.bam # sort bam file\" -> \"sort sample1.bam sample1 # sort bam file\" ---\nconst sortHits = body.search(\"sort sample1.bam -o sample1.bam # sort bam file\", { matchCase: true, matchWholeWord: false });\nsortHits.load(\"text\");\nawait context.sync();\n\nif (sortHits.items.length > 0) {\n  sortHits.items[0].insertText(\"sort sample1.bam sample1 # sort bam file\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Slight change to samtools exercises\n# 1) \"samtools view -bS sample1.sam > sample1.bam\"\n#       -> \"samtools view -S -b sample1.sam > sample1.bam\"\n# 2) \"samtools sort sample1.bam -o sample1.bam # sort bam file\"\n#       -> \"samtools sort sample1.bam sample1 # sort bam file\"\n\n$d = $word.ActiveDocument\n\n# --- Edit 1: \"-bS \" -> \"-S -b \" on the samtools view line ---\n$rng1 = $d.Content\n$rng1.Find.ClearFormatting()\n$rng1.Find.Replacement.ClearFormatting()\n$rng1.Find.Execute(\"-bS \", $false, $false, $false, $false, $false, $true, 1, $false, \"-S -b \", 2)\n\n# --- Edit 2: \"sort sample1.bam -o sample1.bam # sort bam file\" -> \"sort sample1.bam sample1 # sort bam file\" ---\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Replacement.ClearFormatting()\n$rng2.Find.Execute(\"sort sample1.bam -o sample1.bam # sort bam file\", $false, $false, $false, $false, $false, $true, 1, $false, \"sort sample1.bam sample1 # sort bam file\", 2)\n"}
